$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 6.511799999999995
$ws.Range("A3").Value = -21.59790000000002
$ws.Range("B5").Value = 4.808000000000002
$ws.Range("C5").Value = -13.6192
$ws.Range("E7").Value = 11.7916
$ws.Range("C9").Value = -11.61200000000001
$ws.Range("C11").Value = -13.61640000000001
$ws.Range("E11").Value = 13.4474
$ws.Range("A14").Value = -20.38089999999998
$ws.Range("E19").Value = 13.5468
$ws.Range("A21").Value = -21.11510000000001
$ws.Range("C21").Value = -11.358
$ws.Range("E21").Value = 12.387
$ws.Range("A23").Value = -21.47440000000003
$ws.Range("A25").Value = -22.39820000000003
